$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-03 Monday" "2025-03-04 Tuesday"

Replace-Text "27×38=" "40×53="
Replace-Text "97×33=" "30×32="
Replace-Text "72×22=" "48×17="
Replace-Text "34×82=" "59×81="
Replace-Text "76×48=" "63×62="
Replace-Text "90×65=" "29×83="
Replace-Text "85×64=" "12×75="
Replace-Text "66×23=" "12×17="
Replace-Text "91×25=" "16×16="
Replace-Text "49×39=" "46×38="
Replace-Text "84×36=" "94×25="
Replace-Text "27×46=" "89×11="
Replace-Text "44×87=" "92×79="
Replace-Text "60×71=" "21×21="
Replace-Text "57×59=" "28×41="
Replace-Text "33×90=" "82×93="
Replace-Text "96×84=" "71×69="
Replace-Text "46×29=" "55×33="
Replace-Text "49×98=" "61×30="
Replace-Text "50×28=" "88×87="
Replace-Text "29×18=" "69×62="
Replace-Text "80×94=" "80×45="
Replace-Text "61×67=" "27×12="
Replace-Text "15×41=" "25×49="
Replace-Text "15×40=" "19×32="
